$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.106176853179932
$ws.Range("B1").Value = 3.053410530090332
$ws.Range("C1").Value = 6.525791645050049
$ws.Range("D1").Value = 1.854673266410828
$ws.Range("E1").Value = 1.209981083869934
